# Fix wording in the benchmark tables (typo + clearer labels).
# These three cells are reused for all three "iterations = ..." tables
# (rows 10/11/12, 17/18/19, 24/25/26) and drive the embedded charts'
# series names via cell references (Sheet1!$A$10, $A$11, $A$12, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Directly On Machine"
$ws.Range("A11").Value = "Using a Container"
$ws.Range("A12").Value = "Using a VM"

$ws.Range("A17").Value = "Directly On Machine"
$ws.Range("A18").Value = "Using a Container"
$ws.Range("A19").Value = "Using a VM"

$ws.Range("A24").Value = "Directly On Machine"
$ws.Range("A25").Value = "Using a Container"
$ws.Range("A26").Value = "Using a VM"

# Move/resize the third chart ("Chart 5", the iterations=10000 chart) to its
# new anchor position. Values below (in points) are derived from the target
# two-cell anchor (from: col5/333375,row26/66674 - to: col11/600075,row37/4762)
# using this workbook's column width (742156.25 EMU/col) and row height
# (182880 EMU/row, i.e. 14.4pt default row height), converted to points
# (12700 EMU per point).
$chart5 = $ws.ChartObjects().Item(3)
$chart5.Left = 406.484375
$chart5.Top = 379.64992125984253
$chart5.Width = 371.625
$chart5.Height = 153.52503937007873

# Update the selected cell shown when the workbook is reopened.
$ws.Range("D24").Select() | Out-Null
